$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) values per regenerated save_data
$gValues = @{
    2 = 3
    3 = 1
    4 = 2
    5 = 2
    6 = 1
    7 = 2
    8 = 2
    9 = 2
    10 = 2
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    15 = 3
    16 = 2
    17 = 3
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 3
    27 = 1
    28 = 2
    29 = 2
    30 = 0
    31 = 0
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 2
    38 = 1
    39 = 3
    40 = 3
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

